$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'65.763.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "

# Row 3
$ws.Range("D3").Formula = "'2.674.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Formula = "'601.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$ws.Range("D6").Formula = "'157.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Formula = "'0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.52%  "

# Row 9
$ws.Range("E9").Value = "  -0.40%  "

# Row 10
$ws.Range("E10").Value = "  +1.28%  "

# Row 11
$ws.Range("D11").Formula = "'0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").Formula = "'29.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "

# Row 14
$ws.Range("D14").Formula = "'0.0000199"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "

# Row 15
$ws.Range("D15").Formula = "'3.154.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

# Row 16
$ws.Range("D16").Formula = "'65.573.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17
$ws.Range("D17").Formula = "'2.675.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "

# Row 18
$ws.Range("D18").Formula = "'12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("E19").Value = "  -0.48%  "

# Row 20
$ws.Range("E20").Value = "  +2.71%  "

# Row 21
$ws.Range("D21").Formula = "'352.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.32%  "

# Row 22
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("D23").Formula = "'69.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
$ws.Range("E24").Value = "  +6.26%  "

# Row 25
$ws.Range("D25").Formula = "'9.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.77%  "

# Row 26
$ws.Range("E26").Value = "  -4.00%  "

# Row 27
$ws.Range("E27").Value = "  +2.12%  "

# Row 28
$ws.Range("E28").Value = "  -0.94%  "

# Row 29
$ws.Range("E29").Value = "  +0.80%  "

# Row 30
$ws.Range("D30").Formula = "'543.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.77%  "

# Row 31
$ws.Range("E31").Value = "  +0.10%  "

# Row 32
$ws.Range("E32").Value = "  -1.05%  "

# Row 33
$ws.Range("D33").Formula = "'1.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").Formula = "'6.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.26%  "

# Row 35
$ws.Range("E35").Value = "  -0.96%  "

# Row 36
$ws.Range("E36").Value = "  -1.81%  "

# Row 37
$ws.Range("E37").Value = "  -0.78%  "

# Row 38
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("E39").Value = "  -0.95%  "

# Row 40
$ws.Range("D40").Formula = "'158.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.27%  "

# Row 41
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$ws.Range("E42").Value = "  +1.51%  "

# Row 43
$ws.Range("D43").Formula = "'165.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "

# Row 44
$ws.Range("E44").Value = "  -0.87%  "

# Row 45
$ws.Range("D45").Formula = "'0.0618"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.75%  "

# Row 46
$ws.Range("E46").Value = "  -2.26%  "

# Row 47
$ws.Range("D47").Formula = "'23.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "

# Row 48
$ws.Range("E48").Value = "  -0.55%  "

# Row 49
$ws.Range("D49").Formula = "'0.0260"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "

# Row 50
$ws.Range("E50").Value = "  +2.71%  "

# Row 51
$ws.Range("D51").Formula = "'20.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.29%  "
